$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.312990069389343
$ws.Range("B1").Value = 1.287930727005005
$ws.Range("C1").Value = 1.101987242698669
$ws.Range("D1").Value = 1.16438889503479
$ws.Range("E1").Value = 1.007644534111023
